$wb = $excel.ActiveWorkbook

# Data for the new row 97 to append to each of the four sheets.
$rows = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = 45772.95092778935; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; D = "0x01,0x56"; E = "0xd"; F = 400; G = [double]"5.68631262647114e+23"; H = 342; I = 13 },
    @{ Sheet = "ROW35-MID-LIFTER"; A = 45772.80702126157; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; D = "0x01,0x56"; E = "0xe"; F = 400; G = [double]"5.68631262647114e+23"; H = 342; I = 14 },
    @{ Sheet = "ROW02-FE-LIFTER";  A = 45772.95053054398; B = "0x01,0x90"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x01,0x56"; E = "0x3";  F = 400; G = [double]"5.68631262647114e+23"; H = 342; I = 3 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = 45773.0146740625;  B = "0x01,0x90"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x01,0x56"; E = "0x3";  F = 400; G = [double]"9.85046333984776e+23"; H = 342; I = 3 }
)

foreach ($row in $rows) {
    $ws = $wb.Worksheets.Item($row.Sheet)
    $r = 97

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
}
